$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.164663916758826
$ws.Range("B5").Value = 0.164663916758826

$ws.Range("T4").Value = 0.17905006976352
$ws.Range("C21").Value = 0.17905006976352

$ws.Range("P5").Value = 0.20705099127268
$ws.Range("D17").Value = 0.20705099127268

$ws.Range("Q15").Value = 0.164690930189243
$ws.Range("N18").Value = 0.164690930189243
